$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 168 (shifts existing rows 168-228 down to 169-229)
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new weekly price record
$ws.Range("A168").Value2 = 4
$ws.Range("B168").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C168").Value2 = "Los Lagos"
$ws.Range("D168").Value2 = 44553
$ws.Range("D168").NumberFormat = $ws.Range("D169").NumberFormat
$ws.Range("E168").Value2 = 10
$ws.Range("F168").Value2 = 100114014
$ws.Range("G168").Value2 = "Betarraga"
$ws.Range("H168").Value2 = "Sin especificar"
$ws.Range("I168").Value2 = "Primera"
$ws.Range("J168").Value2 = 500
$ws.Range("K168").Value2 = 1000
$ws.Range("L168").Value2 = 1000
$ws.Range("M168").Value2 = 1000
$ws.Range("N168").Value2 = "$/paquete 5 unidades"
$ws.Range("O168").Value2 = "Región del Maule"
$ws.Range("P168").Value2 = 200
$ws.Range("Q168").Value2 = 5
$ws.Range("R168").Value2 = "Hortaliza"
